# Updates the vm_pu.xlsx results sheet for the "380 kV" case run (Case_4_69):
# bus 0 voltage setpoint changed from 1.05 p.u. to 1.02 p.u., and the resulting
# per-bus voltage magnitudes (columns B-F, I-N; G stays 1, H has no data) for
# every time step (rows 2-25) are updated to the new power-flow solution.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B"="1.02"; "C"="1.034054088928839"; "D"="1.041977205839834"; "E"="1.03769802115549"; "F"="1.050932742371011"; "I"="1.03925445930406"; "J"="1.039175437204703"; "K"="1.044755195862259"; "L"="1.040488163928564"; "M"="1.05368564625952"; "N"="1.017020389338579" }
    3 = @{ "B"="1.02"; "C"="1.034857388077937"; "D"="1.042592530624964"; "E"="1.038449155713989"; "F"="1.051685093904593"; "I"="1.039431802627169"; "J"="1.039622447759569"; "K"="1.045181770947775"; "L"="1.041049315052982"; "M"="1.054250698096534"; "N"="1.017170038628497" }
    4 = @{ "B"="1.02"; "C"="1.035377923256704"; "D"="1.042991347387656"; "E"="1.038936268127839"; "F"="1.052172928154165"; "I"="1.039545769527719"; "J"="1.039911762484625"; "K"="1.045457759314049"; "L"="1.041412829913182"; "M"="1.054616669633593"; "N"="1.017266856342422" }
    5 = @{ "B"="1.02"; "C"="1.035596933233146"; "D"="1.043159166186373"; "E"="1.039141306282102"; "F"="1.052378253968077"; "I"="1.039593492260323"; "J"="1.040033405580422"; "K"="1.045573775259401"; "L"="1.041565748726095"; "M"="1.054770605020546"; "N"="1.017307554335936" }
    6 = @{ "B"="1.02"; "C"="1.035633716297484"; "D"="1.043187352807651"; "E"="1.039175748049502"; "F"="1.052412743105692"; "I"="1.039601494014759"; "J"="1.040053830843138"; "K"="1.04559325424654"; "L"="1.041591430118918"; "M"="1.054796456147258"; "N"="1.017314387449125" }
    7 = @{ "B"="1.02"; "C"="1.035380848985573"; "D"="1.042993589179153"; "E"="1.038939006853938"; "F"="1.052175670786182"; "I"="1.03954640794461"; "J"="1.039913387827631"; "K"="1.045459309564085"; "L"="1.041414872842062"; "M"="1.054618726209805"; "N"="1.017267400168018" }
    8 = @{ "B"="1.02"; "C"="1.034325411885392"; "D"="1.042185019839933"; "E"="1.037951646154397"; "F"="1.051186792188061"; "I"="1.039314555616943"; "J"="1.039326491424145"; "K"="1.044899365176989"; "L"="1.040677721192112"; "M"="1.053876535677116"; "N"="1.017070966907136" }
    9 = @{ "B"="1.02"; "C"="1.032471398201878"; "D"="1.040765359362017"; "E"="1.036220136559583"; "F"="1.049452109133636"; "I"="1.0389000154446"; "J"="1.038292892530335"; "K"="1.043912469887122"; "L"="1.039381996390712"; "M"="1.052571419932997"; "N"="1.016724730278137" }
    10 = @{ "B"="1.02"; "C"="1.031239398761336"; "D"="1.039822489216967"; "E"="1.035071528576552"; "F"="1.048301055510246"; "I"="1.038619675828376"; "J"="1.037604301469522"; "K"="1.043254485851431"; "L"="1.038520447753125"; "M"="1.051703272003519"; "N"="1.016493870550663" }
    11 = @{ "B"="1.02"; "C"="1.030706903526488"; "D"="1.039415086796612"; "E"="1.034575553639362"; "F"="1.047803944448563"; "I"="1.038497351456417"; "J"="1.037306263472809"; "K"="1.042969576056908"; "L"="1.038147945899569"; "M"="1.05132783341828"; "N"="1.016393903573368" }
    12 = @{ "B"="1.02"; "C"="1.030509258482954"; "D"="1.039263891544256"; "E"="1.034391535858566"; "F"="1.047619493163157"; "I"="1.038451774965387"; "J"="1.037195579260088"; "K"="1.042863749521797"; "L"="1.038009667065088"; "M"="1.051188451947528"; "N"="1.016356771409789" }
    13 = @{ "B"="1.02"; "C"="1.030551647315729"; "D"="1.039296317427911"; "E"="1.034430998788305"; "F"="1.047659049599899"; "I"="1.038461557590048"; "J"="1.037219320455727"; "K"="1.042886449575718"; "L"="1.038039324474328"; "M"="1.0512183464086"; "N"="1.016364736375832" }
    14 = @{ "B"="1.02"; "C"="1.030690563107891"; "D"="1.039402586242861"; "E"="1.034560338388596"; "F"="1.047788693598178"; "I"="1.038493586936309"; "J"="1.037297113857818"; "K"="1.042960828360771"; "L"="1.038136513985702"; "M"="1.051316310600393"; "N"="1.016390834213774" }
    15 = @{ "B"="1.02"; "C"="1.03077617329954"; "D"="1.039468079536082"; "E"="1.034640056594025"; "F"="1.047868597832512"; "I"="1.038513302775476"; "J"="1.037345047680404"; "K"="1.043006655845396"; "L"="1.038196406955933"; "M"="1.051376679315565"; "N"="1.016406913971109" }
    16 = @{ "B"="1.02"; "C"="1.031274759284037"; "D"="1.039849545598373"; "E"="1.035104474050144"; "F"="1.048334074767243"; "I"="1.038627774450587"; "J"="1.037624084051267"; "K"="1.043273394530344"; "L"="1.038545181295467"; "M"="1.051728198786624"; "N"="1.016500505005356" }
    17 = @{ "B"="1.02"; "C"="1.031587769752509"; "D"="1.040089062567332"; "E"="1.035396161757545"; "F"="1.048626406469119"; "I"="1.038699329600213"; "J"="1.037799151021589"; "K"="1.043440714088233"; "L"="1.0387641078145"; "M"="1.051948826118435"; "N"="1.016559211688724" }
    18 = @{ "B"="1.02"; "C"="1.031770436798653"; "D"="1.040228852244063"; "E"="1.03556643126316"; "F"="1.048797044156125"; "I"="1.038740976197477"; "J"="1.037901276734923"; "K"="1.043538308819243"; "L"="1.03889185732937"; "M"="1.052077560041139"; "N"="1.016593453962263" }
    19 = @{ "B"="1.02"; "C"="1.03183273731543"; "D"="1.040276531011709"; "E"="1.035624511260368"; "F"="1.048855248427244"; "I"="1.038755161269183"; "J"="1.037936100976704"; "K"="1.043571586068574"; "L"="1.038935425616655"; "M"="1.052121462686022"; "N"="1.016605129619403" }
    20 = @{ "B"="1.02"; "C"="1.031554177023946"; "D"="1.04006335601441"; "E"="1.035364852655148"; "F"="1.048595029042097"; "I"="1.038691661745577"; "J"="1.037780366729887"; "K"="1.043422762266466"; "L"="1.038740613537653"; "M"="1.051925150150805"; "N"="1.016552913048907" }
    21 = @{ "B"="1.02"; "C"="1.030649651785627"; "D"="1.039371289061988"; "E"="1.034522245303017"; "F"="1.047750511191618"; "I"="1.038484158944361"; "J"="1.037274205066574"; "K"="1.042938925605657"; "L"="1.038107891739005"; "M"="1.051287460554692"; "N"="1.01638314904631" }
    22 = @{ "B"="1.02"; "C"="1.03008179458204"; "D"="1.038936924622466"; "E"="1.033993677136189"; "F"="1.047220675344026"; "I"="1.03835288539797"; "J"="1.036956079667346"; "K"="1.042634728523366"; "L"="1.037710566681222"; "M"="1.050886943801518"; "N"="1.0162764120603" }
    23 = @{ "B"="1.02"; "C"="1.03038274488348"; "D"="1.03916711621631"; "E"="1.034273765495421"; "F"="1.047501441994914"; "I"="1.038422552348331"; "J"="1.037124712268853"; "K"="1.042795987714418"; "L"="1.037921149002137"; "M"="1.051099224537729"; "N"="1.016332995152445" }
    24 = @{ "B"="1.02"; "C"="1.031569355843852"; "D"="1.040074971441473"; "E"="1.035378999480168"; "F"="1.048609206763105"; "I"="1.038695126798175"; "J"="1.037788854505506"; "K"="1.043430873921986"; "L"="1.038751229425345"; "M"="1.051935848160384"; "N"="1.016555759134241" }
    25 = @{ "B"="1.02"; "C"="1.032950006451178"; "D"="1.041131754619706"; "E"="1.036666771837551"; "F"="1.049899623379791"; "I"="1.039007888779066"; "J"="1.038560024765945"; "K"="1.044167621453236"; "L"="1.039716579592897"; "M"="1.052908491007601"; "N"="1.016814249114637" }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = [double]$data[$row][$col]
    }
}
